# Auto update: 2025-12-03 03:05:24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed data pull - date moves to 2025-12-03 and the row order / figures
# are recomputed for this run.
# Columns: A 날짜, B 종목명, C 티커, D 종가, E RSI, F 5일수익률, G 점수(룰),
#          H 3일상승확률(%), I 5일상승확률(%), J 10일상승확률(%), K 최종점수,
#          L 예측방식, M 판단, N MACRO_SCORE, O MACRO_SIGNAL

$newDate = "2025-12-03"
$macroScore = 65.32892478746797
$macroSignal = "🟢 상승 우위 (다소 완화)"
$judgment = "⛔ 관망하십시오."

$rows = @(
    @{ Name = "Riot Platforms, Inc."; Ticker = "RIOT"; Close = 15.68; RSI = 47.1; Ret5 = 12.93; Rule = 60; P3 = 43; P5 = 46; P10 = 60; Final = 56 },
    @{ Name = "Bitcoin USD"; Ticker = "BTC-USD"; Close = 91512.91; RSI = 47.3; Ret5 = 0.25; Rule = 50; P3 = 30; P5 = 40; P10 = 46; Final = 50.6 },
    @{ Name = "Coinbase Global, Inc."; Ticker = "COIN"; Close = 268.74; RSI = 36.2; Ret5 = 4.99; Rule = 30; P3 = 43; P5 = 53; P10 = 60; Final = 49.8 },
    @{ Name = "MARA Holdings, Inc."; Ticker = "MARA"; Close = 12.15; RSI = 34.8; Ret5 = 8.34; Rule = 30; P3 = 46; P5 = 50; P10 = 60; Final = 48.6 },
    @{ Name = "Strategy Inc"; Ticker = "MSTR"; Close = 183.43; RSI = 30.3; Ret5 = 2.45; Rule = 30; P3 = 40; P5 = 36; P10 = 40; Final = 43 }
)

# Column A holds the date as plain text (not a real date serial) in this
# workbook, so force Text format before writing it - otherwise Excel's
# COM layer auto-converts the "yyyy-mm-dd" literal into a date value.
$dateRange = $ws.Range("A2:A6")
$dateRange.NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 2).Value = $row.Name
    $ws.Cells.Item($r, 3).Value = $row.Ticker
    $ws.Cells.Item($r, 4).Value = $row.Close
    $ws.Cells.Item($r, 5).Value = $row.RSI
    $ws.Cells.Item($r, 6).Value = $row.Ret5
    $ws.Cells.Item($r, 7).Value = $row.Rule
    $ws.Cells.Item($r, 8).Value = $row.P3
    $ws.Cells.Item($r, 9).Value = $row.P5
    $ws.Cells.Item($r, 10).Value = $row.P10
    $ws.Cells.Item($r, 11).Value = $row.Final
    $ws.Cells.Item($r, 12).Value = "Pattern"
    $ws.Cells.Item($r, 13).Value = $judgment
    $ws.Cells.Item($r, 14).Value = $macroScore
    $ws.Cells.Item($r, 15).Value = $macroSignal
    $r++
}

# Restore the default (unformatted) cell style now that the text has been
# committed, so column A doesn't keep a lingering explicit "@" style.
$dateRange.Style = "Normal"
